$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8's "Test scenario" cell (column D) is being updated to a new test case:
# out with "special characters", in with "different languages".
$ws.Range("D8").Value = "Verify the search field handles different languages"

# Move/save the active selection to E3, matching the workbook's saved cursor state.
$ws.Range("E3").Select()
